$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 513, shifting existing rows 513:536 down to 514:536.
$ws.Rows.Item(513).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the new row 513 with its data values.
$ws.Cells.Item(513, 1).Value = 10
$ws.Cells.Item(513, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(513, 3).Value = "La Araucanía"
$ws.Cells.Item(513, 4).Value = 45008
$ws.Cells.Item(513, 5).Value = 9
$ws.Cells.Item(513, 6).Value = 100112040
$ws.Cells.Item(513, 7).Value = "Cilantro"
$ws.Cells.Item(513, 8).Value = "Sin especificar"
$ws.Cells.Item(513, 9).Value = "Primera"
$ws.Cells.Item(513, 10).Value = 65
$ws.Cells.Item(513, 11).Value = 5000
$ws.Cells.Item(513, 12).Value = 5000
$ws.Cells.Item(513, 13).Value = 5000
$ws.Cells.Item(513, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(513, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(513, 16).Value = 2500
$ws.Cells.Item(513, 17).Value = 2
$ws.Cells.Item(513, 18).Value = "Hortaliza"
